# Weekly fruit/vegetable price update ("Fruta / hortaliza, semanal")
# Insert 6 new daily price records (date 45034 = 2023-04-18) for
# Zanahoria / Agrícola del Norte S.A. de Arica, before the existing
# row 427, shifting all subsequent rows down by 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new blank rows at the top of the data block (row 427..432).
# Excel will shift the existing rows 427-498 down to 433-504 and will
# also extend the sheet dimension automatically.
$ws.Rows("427:432").Insert()

# Common, fixed values shared by every new record in this block.
$mercadoId = 1
$mercado   = "Agrícola del Norte S.A. de Arica"
$region    = "Arica y Parinacota"
$fecha     = 45034
$codreg    = 15
$catId     = 100114013
$categoria = "Zanahoria"
$variedad  = "Sin especificar"
$unidad    = "`$/saco 25 kilos"
$kgUnidades = 25
$clasificacion = "Hortaliza"

function Set-PrecioRow {
    param(
        [int]$rowNum,
        [string]$calidad,
        [double]$volumen,
        [double]$precioMin,
        [double]$precioMax,
        [double]$precioProm,
        [string]$origen,
        [double]$precioKg
    )

    $ws.Cells.Item($rowNum, 1).Value  = $mercadoId
    $ws.Cells.Item($rowNum, 2).Value  = $mercado
    $ws.Cells.Item($rowNum, 3).Value  = $region
    $ws.Cells.Item($rowNum, 4).Value  = $fecha
    $ws.Cells.Item($rowNum, 5).Value  = $codreg
    $ws.Cells.Item($rowNum, 6).Value  = $catId
    $ws.Cells.Item($rowNum, 7).Value  = $categoria
    $ws.Cells.Item($rowNum, 8).Value  = $variedad
    $ws.Cells.Item($rowNum, 9).Value  = $calidad
    $ws.Cells.Item($rowNum, 10).Value = $volumen
    $ws.Cells.Item($rowNum, 11).Value = $precioMin
    $ws.Cells.Item($rowNum, 12).Value = $precioMax
    $ws.Cells.Item($rowNum, 13).Value = $precioProm
    $ws.Cells.Item($rowNum, 14).Value = $unidad
    $ws.Cells.Item($rowNum, 15).Value = $origen
    $ws.Cells.Item($rowNum, 16).Value = $precioKg
    $ws.Cells.Item($rowNum, 17).Value = $kgUnidades
    $ws.Cells.Item($rowNum, 18).Value = $clasificacion
}

Set-PrecioRow 427 "Primera" 60 18000 19000 18500 "Región de Arica y Parinacota" 740
Set-PrecioRow 428 "Primera" 40 18000 19000 18375 "Valle de Camiña" 735
Set-PrecioRow 429 "Segunda" 40 14000 15000 14500 "Región de Arica y Parinacota" 580
Set-PrecioRow 430 "Segunda" 35 14000 15000 14429 "Valle de Camiña" 577
Set-PrecioRow 431 "Tercera" 50 10000 11000 10500 "Región de Arica y Parinacota" 420
Set-PrecioRow 432 "Tercera" 25 10000 10000 10000 "Valle de Camiña" 400
